$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.858.75"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "1.830.24"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.07"
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4615"
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3668"
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07161"
$ws.Range("E9").Value = "  -2.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8756"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07897"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.56"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("D13").Value = "1.870.17"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.333"
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.379"
$ws.Range("E15").Value = "  -3.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.73"
$ws.Range("E16").Value = "  -5.17%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008718"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "26.885.59"
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.43"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.997"
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.43"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.987"
$ws.Range("E24").Value = "  +4.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.70"
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.21"
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.964"
$ws.Range("E27").Value = "  -5.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.49"
$ws.Range("E28").Value = "  -2.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.933"
$ws.Range("E29").Value = "  -3.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08840"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.124"
$ws.Range("E31").Value = "  +3.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7538"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.452"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.125"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.569"
$ws.Range("E35").Value = "  -2.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.085"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01934"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.924"
$ws.Range("E38").Value = "  -2.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05125"
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.891"
$ws.Range("E40").Value = "  -3.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4964"
$ws.Range("E41").Value = "  -4.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1595"
$ws.Range("E42").Value = "  -3.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.299"
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4671"
$ws.Range("E44").Value = "  -3.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.005"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.08"
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.22"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.610"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06095"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.55"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.29"
$ws.Range("E51").Value = "  -2.44%  "
